$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refreshed Price (D) and Volume(1h) (E) figures,
# plus the Hedera/VeChain rows (35/36) swapping rank order.
# Numeric-looking Price strings get a leading apostrophe so Excel keeps
# them as text (matching the source data) instead of coercing to numbers.

$ws.Range("D2").Value = "23.250.14"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").Value = "1.605.66"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'1.000"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").Value = "'304.43"
$ws.Range("E6").Value = "  +0.86%  "

$ws.Range("E7").Value = "  -0.48%  "

$ws.Range("D8").Value = "'52.48"
$ws.Range("E8").Value = "  +4.80%  "

$ws.Range("D9").Value = "'0.3626"
$ws.Range("E9").Value = "  -0.49%  "

$ws.Range("E10").Value = "  +1.65%  "

$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").Value = "'22.94"
$ws.Range("E13").Value = "  +2.51%  "

$ws.Range("D14").Value = "'6.600"
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("D15").Value = "'7.368"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").Value = "1.601.02"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").Value = "'94.08"
$ws.Range("E18").Value = "  +2.60%  "

$ws.Range("D19").Value = "'0.06930"
$ws.Range("E19").Value = "  +1.66%  "

$ws.Range("D20").Value = "'18.17"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").Value = "'6.543"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").Value = "'12.92"
$ws.Range("E23").Value = "  -0.64%  "

$ws.Range("D24").Value = "23.240.79"
$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("D25").Value = "'2.457"
$ws.Range("E25").Value = "  +4.00%  "

$ws.Range("D26").Value = "'3.074"
$ws.Range("E26").Value = "  +9.45%  "

$ws.Range("D27").Value = "'21.21"
$ws.Range("E27").Value = "  +0.67%  "

$ws.Range("D28").Value = "'150.04"
$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("D29").Value = "'5.271"
$ws.Range("E29").Value = "  +0.54%  "

$ws.Range("D30").Value = "'135.16"
$ws.Range("E30").Value = "  +0.66%  "

$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("D32").Value = "'6.730"
$ws.Range("E32").Value = "  -2.61%  "

$ws.Range("D33").Value = "1.776.51"
$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("D34").Value = "'0.9625"
$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.02771"
$ws.Range("E35").Value = "  +2.60%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.07467"
$ws.Range("E36").Value = "  -1.49%  "

$ws.Range("D37").Value = "'10.35"
$ws.Range("E37").Value = "  +0.90%  "

$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("D39").Value = "'6.122"
$ws.Range("E39").Value = "  -1.67%  "

$ws.Range("E40").Value = "  -0.30%  "

$ws.Range("D41").Value = "'1.409"
$ws.Range("E41").Value = "  +3.03%  "

$ws.Range("D42").Value = "'0.7094"
$ws.Range("E42").Value = "  +1.02%  "

$ws.Range("D43").Value = "'12.45"
$ws.Range("E43").Value = "  +0.58%  "

$ws.Range("D44").Value = "'15.84"
$ws.Range("E44").Value = "  +4.23%  "

$ws.Range("D45").Value = "'0.6544"
$ws.Range("E45").Value = "  -1.09%  "

$ws.Range("D46").Value = "'2.333"
$ws.Range("E46").Value = "  +2.04%  "

$ws.Range("D47").Value = "'4.010"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").Value = "'134.08"
$ws.Range("E48").Value = "  +2.02%  "

$ws.Range("D49").Value = "'0.07953"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").Value = "'1.204"
$ws.Range("E50").Value = "  -1.36%  "

$ws.Range("D51").Value = "'1.189"
$ws.Range("E51").Value = "  -3.32%  "
